# Insert a new data row at row 127 (pushing the existing row 127..231 down
# to 128..232) and populate it with the new weekly price-record values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 127; everything below shifts down by one,
# which is what grows the used range from A1:R231 to A1:R232.
$ws.Rows.Item(127).Insert()

$ws.Range("A127").Value = 9
$ws.Range("B127").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C127").Value = "Metropolitana"
$ws.Range("D127").Value = 44596
$ws.Range("E127").Value = 13
$ws.Range("F127").Value = 100112043
$ws.Range("G127").Value = "Pepino ensalada"
$ws.Range("H127").Value = "Sin especificar"
$ws.Range("I127").Value = "Primera"
$ws.Range("J127").Value = 270
$ws.Range("K127").Value = 13000
$ws.Range("L127").Value = 14000
$ws.Range("M127").Value = 13333
$ws.Range("N127").Value = "`$/caja 70 unidades"
$ws.Range("O127").Value = "Región Metropolitana"
$ws.Range("P127").Value = 190
$ws.Range("Q127").Value = 70
$ws.Range("R127").Value = "Hortaliza"
